# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'61.371.77"
$ws.Range('D2').Cells.Item(1,1).Style = "Normal"
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = "'2.927.61"
$ws.Range('D3').Cells.Item(1,1).Style = "Normal"
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Cells.Item(1,1).Style = "Normal"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'597.46"
$ws.Range('D5').Cells.Item(1,1).Style = "Normal"
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').Value = "'145.13"
$ws.Range('D6').Cells.Item(1,1).Style = "Normal"
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('E14').Value = '  +0.19%  '
$ws.Range('D15').Value = "'3.412.91"
$ws.Range('D15').Cells.Item(1,1).Style = "Normal"
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').Value = "'61.349.80"
$ws.Range('D16').Cells.Item(1,1).Style = "Normal"
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').Value = "'2.929.94"
$ws.Range('D17').Cells.Item(1,1).Style = "Normal"
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').Value = "'431.67"
$ws.Range('D19').Cells.Item(1,1).Style = "Normal"
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').Value = "'13.49"
$ws.Range('D20').Cells.Item(1,1).Style = "Normal"
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').Value = "'81.78"
$ws.Range('D23').Cells.Item(1,1).Style = "Normal"
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('D24').Value = "'10.86"
$ws.Range('D24').Cells.Item(1,1).Style = "Normal"
$ws.Range('E24').Value = '  -1.99%  '
$ws.Range('E25').Value = '  -1.69%  '
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('D28').Value = "'2.21"
$ws.Range('D28').Cells.Item(1,1).Style = "Normal"
$ws.Range('E28').Value = '  -4.60%  '
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('D30').Value = "'6.90"
$ws.Range('D30').Cells.Item(1,1).Style = "Normal"
$ws.Range('E30').Value = '  -2.67%  '
$ws.Range('D31').Value = "'26.58"
$ws.Range('D31').Cells.Item(1,1).Style = "Normal"
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  +1.30%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').Value = "'0.0₃0885"
$ws.Range('D34').Cells.Item(1,1).Style = "Normal"
$ws.Range('E34').Value = '  +3.61%  '
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('E37').Value = '  -1.70%  '
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('E39').Value = '  -1.67%  '
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('D41').Value = "'42.04"
$ws.Range('D41').Cells.Item(1,1).Style = "Normal"
$ws.Range('E41').Value = '  +5.23%  '
$ws.Range('E42').Value = '  -2.23%  '
$ws.Range('E43').Value = '  -0.44%  '
$ws.Range('D44').Value = "'2.699.61"
$ws.Range('D44').Cells.Item(1,1).Style = "Normal"
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('D45').Value = "'133.63"
$ws.Range('D45').Cells.Item(1,1).Style = "Normal"
$ws.Range('E45').Value = '  +2.50%  '
$ws.Range('D46').Value = "'364.04"
$ws.Range('D46').Cells.Item(1,1).Style = "Normal"
$ws.Range('E46').Value = '  -3.19%  '
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = "'23.61"
$ws.Range('D48').Cells.Item(1,1).Style = "Normal"
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('E49').Value = '  -1.11%  '
$ws.Range('D50').Value = "'2.00"
$ws.Range('D50').Cells.Item(1,1).Style = "Normal"
$ws.Range('E50').Value = '  -1.41%  '
$ws.Range('E51').Value = '  -2.16%  '
